$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.940.14'
$ws.Range('E2').Value = '  -0.76%  '
$ws.Range('D3').Value = '2.319.23'
$ws.Range('E3').Value = '  -1.68%  '
$ws.Range('E4').Value = '  -0.02%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '530.13'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +1.78%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '132.45'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -2.36%  '
$ws.Range('E7').Value = '  -0.25%  '
$ws.Range('E8').Value = '  -1.36%  '
$ws.Range('D9').Value = '2.343.04'
$ws.Range('E9').Value = '  -1.44%  '
$ws.Range('E10').Value = '  -1.19%  '
$ws.Range('E11').Value = '  -0.12%  '
$ws.Range('E12').Value = '  -2.92%  '
$ws.Range('E13').Value = '  +1.18%  '
$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '23.50'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -3.68%  '
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '2.737.08'
$ws.Range('E15').Value = '  -1.65%  '
$ws.Range('D16').Value = '56.968.44'
$ws.Range('E16').Value = '  -0.71%  '
$ws.Range('E17').Value = '  -2.18%  '
$ws.Range('D18').Value = '2.335.14'
$ws.Range('E18').Value = '  -0.40%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '337.13'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +2.22%  '
$ws.Range('E20').Value = '  -1.81%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '6.87'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +2.01%  '
$ws.Range('E22').Value = '  -1.87%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +0.07%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '61.56'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +0.50%  '
$ws.Range('E25').Value = '  +0.84%  '
$ws.Range('E26').Value = '  -0.14%  '
$ws.Range('E27').Value = '  -0.16%  '
$ws.Range('E28').Value = '  +1.33%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '173.00'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +3.13%  '
$ws.Range('E30').Value = '  +1.03%  '
$ws.Range('D31').Value = '0.0₃0725'
$ws.Range('E31').Value = '  -2.65%  '
$ws.Range('E32').Value = '  -2.71%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '18.50'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -0.40%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '0.998'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -0.05%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '0.992'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  -0.21%  '
$ws.Range('E36').Value = '  -3.28%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.927'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -0.16%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '3.98'
$c.Style = 'Normal'
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '39.25'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +1.10%  '
$ws.Range('E40').Value = '  -3.38%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '5.77'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +7.77%  '
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '149.03'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -0.53%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '0.375'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -2.58%  '
$ws.Range('E44').Value = '  -1.48%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '282.69'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -0.54%  '
$ws.Range('E46').Value = '  -1.41%  '
$ws.Range('E47').Value = '  -1.65%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '18.86'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +3.46%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.559'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -1.12%  '
$ws.Range('E50').Value = '  -1.24%  '
$ws.Range('E51').Value = '  -0.90%  '
